# Scheduled data refresh: update coin price / 1h-volume figures, and
# re-sync a handful of rows whose ranking order shifted (symbol list
# refresh, GitHub Actions scheduled run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'323.38"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-2.67%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'44.46"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'0.74%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.503"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-4.60%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08030"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-3.72%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'8.671"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-1.52%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'4.337"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-3.76%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.884"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-3.96%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.691"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-6.99%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.9414"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'1.00%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1176"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-5.47%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.1876"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-3.38%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.09913"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'4.42%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.04165"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'5.22%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.1064"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.09%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.001270"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-2.44%"
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "'0.005884"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.73%"
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.602"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'2.78%"
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "'0.3485"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-0.69%"
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value = "MCDex"
$ws.Range("C20").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D20").Value = "'8.476"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-6.59%"
$ws.Range("E20").Style = "Normal"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").Value = "'0.1375"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'0.93%"
$ws.Range("E21").Style = "Normal"
$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D22").Value = "'0.2538"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-1.33%"
$ws.Range("E22").Style = "Normal"
$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D23").Value = "'0.04253"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-3.70%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001243"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-1.08%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004467"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'1.84%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001237"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'3.87%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0004008"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'0.39%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.02629"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-7.20%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.05485"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-5.07%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007730"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-2.44%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1389"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-2.59%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.006790"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-25.28%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.001989"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-5.41%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.009223"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-11.93%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00007136"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-1.90%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000754"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'0.37%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.003403"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'5.69%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.002281"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'0.05%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002110"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.37%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002010"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.37%"
$ws.Range("E51").Style = "Normal"
